$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (column A and B got narrower) ---
$ws.Columns.Item(1).ColumnWidth = 38.666666666666664
$ws.Columns.Item(2).ColumnWidth = 37.166666666666664

# --- Add a new "2022" data column (S) mirroring the existing year columns ---

# S4: year header 2022 (copy formatting from R4)
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# S5: Education row value
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 44

# S6: Health row value (uses a "0.0" number format)
$ws.Range("O6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("S6").Value = 20.6
$ws.Range("S6").NumberFormat = "0.0"

# S7: Social protection row value
$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("S7").Value = 7.9

# S8: total row value
$ws.Range("R8").Copy()
$ws.Range("S8").PasteSpecial(-4122)
$ws.Range("S8").Value = 15.5

# --- Restore the active selection like in the authored workbook ---
$ws.Range("Y14").Select()
